$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 156, shifting existing rows 156:222 down to 157:223.
$ws.Rows.Item(156).Insert()

# Populate the newly inserted row 156 with the new weekly price record.
# Columns A,B,C,E,F,G,H,N,O,Q,R keep the same values as the rest of this
# market/category block; D,I,J,K,L,M,P carry the new data point.
$ws.Cells.Item(156, 1).Value = 7
$ws.Cells.Item(156, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(156, 3).Value = "Ñuble"
$ws.Cells.Item(156, 4).Value = 44510
$ws.Cells.Item(156, 5).Value = 16
$ws.Cells.Item(156, 6).Value = 100114001
$ws.Cells.Item(156, 7).Value = "Papa"
$ws.Cells.Item(156, 8).Value = "Patagonia"
$ws.Cells.Item(156, 9).Value = "1a (guarda)"
$ws.Cells.Item(156, 10).Value = 200
$ws.Cells.Item(156, 11).Value = 7000
$ws.Cells.Item(156, 12).Value = 8000
$ws.Cells.Item(156, 13).Value = 7500
$ws.Cells.Item(156, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(156, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(156, 16).Value = 300
$ws.Cells.Item(156, 17).Value = 25
$ws.Cells.Item(156, 18).Value = "Hortaliza"
